# "Add pictures Dec 2"
#
# The real content edit buried in this diff is tiny: a handful of
# scientific-name cells in column B used a doubled space (e.g. before
# "spp." / "commune" / the en-dash) and those doubled spaces were
# collapsed down to a single space. Everything else in the raw XML diff
# (shared-string index churn, sheetView selection/scroll state) is a
# mechanical side effect of Excel rewriting the shared-strings table and
# recording the current selection when it saved - so we reproduce the
# same user-visible actions here: fix the text, then select B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> corrected column-B text (doubled spaces collapsed to one).
$fixes = @{
    3   = "Dicranum spp."
    5   = "Mnium spp."
    7   = "Polytrichum commune"
    9   = "Sphagnum spp."
    12  = "Dryopteris spp."
    13  = "Equisetum spp."
    39  = "Carex spp. $([char]0x2013) family only"
    53  = "Goodyera spp."
    55  = "Juncus spp. $([char]0x2013) family only"
    58  = "Lupinus spp."
    70  = "Poaceae $([char]0x2013) family only"
    74  = "Solidago spp."
    81  = "Vicia spp."
    82  = "Viola spp."
    86  = "Usnea spp."
    90  = "Amelanchier spp."
    98  = "Crataegus spp."
    119 = "Rosa spp."
    122 = "Salix spp."
}

foreach ($row in $fixes.Keys) {
    $ws.Cells.Item($row, 2).Value = $fixes[$row]
}

# Match the workbook's recorded selection after the edit (B3, no frozen
# scroll position left over from a prior view).
$ws.Range("B3").Select() | Out-Null
